# Refresh the cryptocurrency price / 1h-volume figures in the sheet,
# matching the "Updated cryptos list on Sat Aug 12 12:51:23 UTC 2023 with
# GitHub Actions" data refresh. This also covers the two data rows (50/51)
# where "RenderToken" and "EnergySwap" swapped ranking positions.
#
# Note: the source feed stores numeric-looking figures (prices such as
# "240.65", "1.001", or the thousands/decimal mixed "29.412.92") as plain
# text. Assigning such a string straight to Range.Value lets Excel's COM
# automation silently reinterpret it as a real number (e.g. "1.001" turns
# into the Double 1.0009999999999999, and "240.65" becomes
# 240.65000000000001), which is not what the workbook contains. So for
# every replacement value that looks like a plain number we temporarily
# force the cell to Text format, assign the literal string, and then
# restore the cell's style so no stray formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "29.412.92"
    "E2" = "  +0.03%  "
    "D3" = "1.850.06"
    "E3" = "  +0.08%  "
    "E4" = "  +0.06%  "
    "D5" = "240.65"
    "E5" = "  +0.04%  "
    "D6" = "0.6293"
    "E6" = "  -0.17%  "
    "D7" = "1.001"
    "E7" = "  +0.05%  "
    "D8" = "0.07668"
    "E8" = "  +1.44%  "
    "D9" = "0.2936"
    "E9" = "  -0.67%  "
    "D10" = "24.53"
    "E10" = "  +0.52%  "
    "D11" = "0.07749"
    "E11" = "  +0.68%  "
    "D12" = "1.852.51"
    "E12" = "  -0.04%  "
    "E13" = "  +0.58%  "
    "D14" = "0.00001092"
    "E14" = "  +9.05%  "
    "D15" = "0.6813"
    "E15" = "  -0.45%  "
    "D16" = "83.64"
    "E16" = "  +0.66%  "
    "D17" = "2.104.32"
    "E17" = "  +0.53%  "
    "D18" = "6.161"
    "E18" = "  +0.34%  "
    "D19" = "29.452.83"
    "E19" = "  +0.07%  "
    "D20" = "229.04"
    "E20" = "  +0.64%  "
    "E21" = "  +0.09%  "
    "D22" = "1.001"
    "E22" = "  +0.03%  "
    "E23" = "  -1.29%  "
    "E24" = "  +0.00%  "
    "D25" = "157.14"
    "E25" = "  -0.02%  "
    "D26" = "0.1389"
    "E26" = "  -0.45%  "
    "D27" = "8.404"
    "E27" = "  +0.26%  "
    "E28" = "  -0.09%  "
    "D29" = "1.314"
    "E29" = "  +4.22%  "
    "D30" = "1.465"
    "E30" = "  -0.41%  "
    "D31" = "0.05723"
    "E31" = "  +0.36%  "
    "D32" = "4.126"
    "E32" = "  +0.03%  "
    "E33" = "  +0.90%  "
    "D34" = "1.851"
    "E35" = "  +0.59%  "
    "D36" = "0.7082"
    "E36" = "  -1.04%  "
    "E37" = "  -0.01%  "
    "E38" = "  -0.19%  "
    "D39" = "0.01793"
    "E39" = "  -0.96%  "
    "D40" = "1.217.39"
    "E40" = "  -2.76%  "
    "D41" = "6.482"
    "E41" = "  +4.81%  "
    "D42" = "0.9103"
    "E42" = "  +0.44%  "
    "D43" = "1.001"
    "E43" = "  +0.05%  "
    "D44" = "2.013.14"
    "E44" = "  +0.51%  "
    "D45" = "101.73"
    "E45" = "  +0.47%  "
    "D46" = "66.26"
    "E46" = "  +0.23%  "
    "D47" = "0.00000000120"
    "E47" = "  +2.18%  "
    "D48" = "7.144"
    "E48" = "  +0.52%  "
    "E49" = "  -0.07%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "8.987"
    "E50" = "  -1.24%  "
    "B51" = "RenderToken"
    "C51" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D51" = "1.684"
    "E51" = "  +0.16%  "
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $range = $ws.Range($cellRef)

    if ($newValue.Trim() -match '^[+-]?\d+(\.\d+)?$') {
        # Plain-number-looking text (e.g. "240.65", "1.001") - keep it as
        # text so Excel does not coerce it into a floating point Double.
        $range.NumberFormat = "@"
        $range.Value = $newValue
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}
